$wb = $excel.ActiveWorkbook

# Configurable zero_before_threshold parameter: recomputed First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E), and Pulse_Width (G) on each Step3_DataPts_* sheet
# so that dims before the noise_threshold / First Rise Point are zeroed out.

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.001424441565985998
$ws.Range("G2").Value = 7
$ws.Range("C3").Value = 88
$ws.Range("E3").Value = 0.001372913725065856
$ws.Range("G3").Value = 7
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.001350453340344121
$ws.Range("G4").Value = 12
$ws.Range("C5").Value = 87
$ws.Range("E5").Value = 0.001884271757899121
$ws.Range("G5").Value = 8
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.001077945496345049
$ws.Range("G6").Value = 8

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.001424441565985998
$ws.Range("G2").Value = 27
$ws.Range("C3").Value = 88
$ws.Range("E3").Value = 0.001372913725065856
$ws.Range("G3").Value = 27
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.001350453340344121
$ws.Range("G4").Value = 27
$ws.Range("C5").Value = 87
$ws.Range("E5").Value = 0.001884271757899121
$ws.Range("G5").Value = 28
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.001077945496345049
$ws.Range("G6").Value = 28

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.001424441565985998
$ws.Range("G2").Value = 37
$ws.Range("C3").Value = 88
$ws.Range("E3").Value = 0.001372913725065856
$ws.Range("G3").Value = 41
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.001350453340344121
$ws.Range("G4").Value = 38
$ws.Range("C5").Value = 87
$ws.Range("E5").Value = 0.001884271757899121
$ws.Range("G5").Value = 45
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.001077945496345049
$ws.Range("G6").Value = 46

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.001424441565985998
$ws.Range("G2").Value = 69
$ws.Range("C3").Value = 88
$ws.Range("E3").Value = 0.001372913725065856
$ws.Range("G3").Value = 70
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.001350453340344121
$ws.Range("G4").Value = 70
$ws.Range("C5").Value = 87
$ws.Range("E5").Value = 0.001884271757899121
$ws.Range("G5").Value = 73
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.001077945496345049
$ws.Range("G6").Value = 70
